$d = $word.ActiveDocument

# 1. Replace "APIARY AUTHORITY" heading text with the template placeholder "{{ }}"
#    (MatchCase=$true so the later, lower-case "apiary authority" elsewhere in the
#    document is left untouched)
$d.Content.Find.Execute("APIARY AUTHORITY", $true, $false, $false, $false, $false,
                         $true, 1, $false, "{{ }}", 2)

# 2. Remove the paragraph containing the "CONSERVATION AND LAND MANAGEMENT REGULATIONS 2002 (PART 8A)"
#    text entirely (including its paragraph mark), leaving the following empty paragraph intact.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*CONSERVATION AND LAND MANAGEMENT REGULATIONS 2002 (PART 8A)*") {
        $p.Range.Delete()
        break
    }
}
